# Adicionados novos tratamentos de exceções
# Append a new "agendamento" (appointment) row with the client's name, date,
# time and phone number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The phone number is long enough that Excel would otherwise infer it as a
# number (and lose precision/leading context); format the cell as Text first
# so it is stored/treated as a literal string, same as the other new cells.
$ws.Range("D2").NumberFormat = "@"

$ws.Range("A2").Value = "marcos"
$ws.Range("B2").Value = "21/02/2024"
$ws.Range("C2").Value = "10:00"
$ws.Range("D2").Value = "5533991965662"

# Mirror the author's resulting selection: the block A2:D5 selected (the
# active cell lands on the top-left corner of the selected block, A2, which
# is how Range.Select always anchors a multi-cell selection).
$ws.Range("A2:D5").Select()
